$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds date serial numbers that were updated from
# 45185 (2023-09-16) to 45204 (2023-10-05) for rows 2 through 15.
foreach ($row in 2..15) {
    $ws.Cells.Item($row, 3).Value = 45204
}
